$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.460.16'
$ws.Range("E2").Value = '  -3.23%  '
$ws.Range("D3").Value = '3.155.52'
$ws.Range("E3").Value = '  -2.69%  '
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '608.29'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.62%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.36'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.69%  '
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("D8").Value = '3.146.71'
$ws.Range("E8").Value = '  -2.90%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.527'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.64%  '
$ws.Range("E10").Value = '  -7.54%  '
$ws.Range("E11").Value = '  -3.77%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.474'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.60%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000255'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -6.14%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.00'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -7.79%  '
$ws.Range("D15").Value = '3.670.98'
$ws.Range("E15").Value = '  -2.61%  '
$ws.Range("D16").Value = '64.400.02'
$ws.Range("E16").Value = '  -3.34%  '
$ws.Range("E17").Value = '  +0.97%  '
$ws.Range("D18").Value = '3.152.24'
$ws.Range("E18").Value = '  -2.65%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.93'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '478.24'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.99%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.58'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.76%  '
$ws.Range("E22").Value = '  -5.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.73'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.64%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.75'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -6.23%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.41'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.09%  '
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("E27").Value = '  -3.85%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.40'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -7.67%  '
$ws.Range("E29").Value = '  -6.96%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.116'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -32.25%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.79'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.60%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.75'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.63%  '
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.15'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.09%  '
$ws.Range("E35").Value = '  -4.83%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.01'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '54.14'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.16%  '
$ws.Range("D38").Value = '0.0₃0718'
$ws.Range("E38").Value = '  -11.13%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '451.73'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -8.94%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.91'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -10.22%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0396'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.43'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.48%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.118'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.36%  '
$ws.Range("D44").Value = '2.839.93'
$ws.Range("E44").Value = '  -3.77%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.268'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -8.26%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.26'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -8.50%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '26.39'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.32%  '
$ws.Range("E48").Value = '  -0.08%  '
$ws.Range("E49").Value = '  -4.03%  '
$ws.Range("E50").Value = '  -4.30%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '118.45'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.04%  '
